# Updates cryptocurrency price/volume data to latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.432.63'
$ws.Range('E2').Value = '  +0.81%  '
$ws.Range('D3').Value = '2.008.91'
$ws.Range('E3').Value = '  -0.34%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '260.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.615'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.64%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '56.08'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -6.19%  '
$ws.Range('E9').Value = '  -1.22%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0771'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.49%  '
$ws.Range('E11').Value = '  -2.66%  '
$ws.Range('B12').Value = 'Chainlink'
$ws.Range('C12').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.26'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.82%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '2.303.50'
$ws.Range('E13').Value = '  -0.36%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.03'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.77%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.800'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.75%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.24'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.03%  '
$ws.Range('D17').Value = '1.993.67'
$ws.Range('E17').Value = '  -1.30%  '
$ws.Range('D18').Value = '37.293.53'
$ws.Range('E18').Value = '  +0.73%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.77'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.96%  '
$ws.Range('D20').Value = '0.0₃0837'
$ws.Range('E20').Value = '  -3.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.14'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '230.57'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.61'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.32%  '
$ws.Range('E24').Value = '  -0.06%  '
$ws.Range('E25').Value = '  -1.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.41'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.51%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.89'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.76%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.63'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.53%  '
$ws.Range('E29').Value = '  -4.53%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.34'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.61'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.16%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0645'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.52'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.56%  '
$ws.Range('E35').Value = '  -3.43%  '
$ws.Range('E36').Value = '  +0.46%  '
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.34'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.57%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.37'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.73%  '
$ws.Range('E40').Value = '  +3.87%  '
$ws.Range('E41').Value = '  +0.45%  '
$ws.Range('E42').Value = '  -1.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0924'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.38%  '
$ws.Range('D44').Value = '1.412.70'
$ws.Range('E44').Value = '  +1.98%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '15.78'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '89.38'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.83%  '
$ws.Range('E47').Value = '  -3.33%  '
$ws.Range('E48').Value = '  +2.53%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.02'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.48%  '
$ws.Range('E50').Value = '  -0.30%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.94'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -8.53%  '
